$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "60.902.74"
$ws.Range("E2").Value = "  +0.21%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.367.14"
$ws.Range("E3").Value = "  -0.66%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5: BNB
$ws.Range("D5").Value = "570.98"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6: Solana
$ws.Range("D6").Value = "138.82"
$ws.Range("E6").Value = "  -2.24%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8: XRP
$ws.Range("D8").Value = "0.471"
$ws.Range("E8").Value = "  -0.52%  "

# Row 9: Toncoin
$ws.Range("D9").Value = "7.63"
$ws.Range("E9").Value = "  +1.61%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -2.14%  "

# Row 11: Cardano
$ws.Range("D11").Value = "'0.380"
$ws.Range("E11").Value = "  -3.40%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "3.943.06"
$ws.Range("E12").Value = "  -0.65%  "

# Row 13: TRON
$ws.Range("E13").Value = "  +1.91%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "27.56"
$ws.Range("E14").Value = "  -2.36%  "

# Row 15: WrappedEther
$ws.Range("D15").Value = "3.367.77"
$ws.Range("E15").Value = "  -0.61%  "

# Row 16: ShibaInu
$ws.Range("D16").Value = "0.0000167"
$ws.Range("E16").Value = "  -2.80%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "60.983.20"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18: Polkadot
$ws.Range("D18").Value = "6.04"
$ws.Range("E18").Value = "  -3.74%  "

# Row 19: Chainlink
$ws.Range("E19").Value = "  -3.54%  "

# Row 20: Uniswap
$ws.Range("D20").Value = "8.81"
$ws.Range("E20").Value = "  -2.43%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "380.83"
$ws.Range("E21").Value = "  -1.79%  "

# Row 22: Litecoin
$ws.Range("D22").Value = "75.06"
$ws.Range("E22").Value = "  +1.92%  "

# Row 23: Polygon
$ws.Range("D23").Value = "0.547"
$ws.Range("E23").Value = "  -2.49%  "

# Row 24: Dai
$ws.Range("E24").Value = "  -0.06%  "

# Row 25: PEPE
$ws.Range("D25").Value = "0.0000111"
$ws.Range("E25").Value = "  -5.56%  "

# Row 26: Kaspa
$ws.Range("D26").Value = "0.188"
$ws.Range("E26").Value = "  +5.37%  "

# Row 27: Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.02%  "

# Row 28: RenderToken
$ws.Range("D28").Value = "7.12"
$ws.Range("E28").Value = "  -3.72%  "

# Row 29: InternetComputer(DFINITY)
$ws.Range("D29").Value = "7.83"
$ws.Range("E29").Value = "  -1.96%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  -1.82%  "

# Row 31: USDe
$ws.Range("E31").Value = "  -0.03%  "

# Row 32: Fetch.AI
$ws.Range("D32").Value = "1.34"
$ws.Range("E32").Value = "  -4.94%  "

# Row 33: EthereumClassic
$ws.Range("D33").Value = "22.78"
$ws.Range("E33").Value = "  -3.85%  "

# Row 34: Aptos
$ws.Range("D34").Value = "6.84"
$ws.Range("E34").Value = "  -1.68%  "

# Row 35: Monero
$ws.Range("D35").Value = "166.04"
$ws.Range("E35").Value = "  -0.54%  "

# Row 36: NEARProtocol
$ws.Range("D36").Value = "'4.90"
$ws.Range("E36").Value = "  -1.33%  "

# Row 37: RenzoRestakedETH
$ws.Range("D37").Value = "3.406.14"
$ws.Range("E37").Value = "  -0.38%  "

# Row 38: ImmutableX
$ws.Range("E38").Value = "  -3.79%  "

# Row 39: Hedera
$ws.Range("D39").Value = "0.0757"
$ws.Range("E39").Value = "  -2.52%  "

# Row 40: EnergySwap
$ws.Range("D40").Value = "25.28"
$ws.Range("E40").Value = "  -9.72%  "

# Row 41: Mantle
$ws.Range("D41").Value = "0.771"
$ws.Range("E41").Value = "  -1.43%  "

# Row 42: Filecoin
$ws.Range("D42").Value = "4.32"
$ws.Range("E42").Value = "  -2.49%  "

# Row 43: Stacks
$ws.Range("E43").Value = "  -3.87%  "

# Row 44: ONDO
$ws.Range("D44").Value = "'1.10"
$ws.Range("E44").Value = "  -2.24%  "

# Row 45: Maker
$ws.Range("D45").Value = "2.444.57"
$ws.Range("E45").Value = "  -4.38%  "

# Row 46: FirstDigitalUSD
$ws.Range("E46").Value = "  +0.03%  "

# Row 47: Cosmos
$ws.Range("D47").Value = "6.57"
$ws.Range("E47").Value = "  -3.83%  "

# Row 48: InjectiveProtocol
$ws.Range("D48").Value = "22.17"
$ws.Range("E48").Value = "  -4.72%  "

# Row 49: VeChain
$ws.Range("D49").Value = "0.0258"
$ws.Range("E49").Value = "  -4.43%  "

# Row 50: dogwifhat
$ws.Range("D50").Value = "2.02"
$ws.Range("E50").Value = "  -2.42%  "

# Row 51: TheGraph
$ws.Range("D51").Value = "0.201"
$ws.Range("E51").Value = "  -3.38%  "
